# Update gitbook 2024-09-04 17:35:20
#
# 1. Refresh the "datetimeFigureOut" date placeholder cached text on the two
#    slide layouts that carry it (Agenda / Title-and-Content-3) from
#    8/21/2024 -> 9/4/2024.
# 2. Swap the three replit.com demo links (now hosted on jsfiddle.net) on
#    the Lists / Inputs / Attributes example slides.

$p = $ppt.ActivePresentation

# --- 1. Date placeholders on the slide layouts -----------------------------
$design = $p.Designs.Item(1)
$master = $design.SlideMaster

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    for ($si = 1; $si -le $layout.Shapes.Count; $si++) {
        $shape = $layout.Shapes.Item($si)
        if ($shape.HasTextFrame) {
            $tr = $shape.TextFrame.TextRange
            if ($tr.Text -eq "8/21/2024") {
                $tr.Text = "9/4/2024"
            }
        }
    }
}

# --- 2. replit.com -> jsfiddle.net links on the example slides -------------
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTextFrame) {
            $tr = $shape.TextFrame.TextRange
            switch ($tr.Text) {
                "https://replit.com/@HylandOutreach/ListExample" {
                    $tr.Text = "https://jsfiddle.net/5fmw3dL8/"
                }
                "https://replit.com/@HylandOutreach/InputExamples" {
                    $tr.Text = "https://jsfiddle.net/mz7gpx6c/"
                }
                "https://replit.com/@HylandOutreach/AttributesExample" {
                    $tr.Text = "https://jsfiddle.net/Lhodz460/"
                }
            }
        }
    }
}
